# Auto-generated edit script: restores correct batch (B/C/D/E/F/G) ordering
# for grouped product rows that had been shuffled within their row-blocks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 149
$ws.Range("B149").Value = 48654
$ws.Range("C149").Value = "CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms"
$ws.Range("D149").Value = 32.02
$ws.Range("E149").Value = 38.26
$ws.Range("F149").Value = -1
$ws.Range("G149").Value = -32.02

# Row 150
$ws.Range("B150").Value = 63902
$ws.Range("C150").Value = "CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms"
$ws.Range("D150").Value = 32.02
$ws.Range("E150").Value = 34.04
$ws.Range("F150").Value = 2
$ws.Range("G150").Value = 64.04000000000001

# Row 264
$ws.Range("B264").Value = 48719
$ws.Range("C264").Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("D264").Value = 295.75
$ws.Range("E264").Value = 353.35
$ws.Range("F264").Value = -81
$ws.Range("G264").Value = -23955.75

# Row 265
$ws.Range("B265").Value = 64979
$ws.Range("C265").Value = "HIM-BABY CARE GIFT PACK (WW)1"
$ws.Range("D265").Value = 295.75
$ws.Range("E265").Value = 314.41
$ws.Range("F265").Value = 82
$ws.Range("G265").Value = 24251.5

# Row 316
$ws.Range("B316").Value = 57077
$ws.Range("C316").Value = "HUL-Bru Inst Poly 50g"
$ws.Range("D316").Value = 93.08
$ws.Range("E316").Value = 111.2
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 93.08

# Row 317
$ws.Range("B317").Value = 61610
$ws.Range("C317").Value = "HUL-Bru Inst Poly 50g"
$ws.Range("D317").Value = 102.71
$ws.Range("E317").Value = 122.71
$ws.Range("F317").Value = -58
$ws.Range("G317").Value = -5957.18

# Row 318
$ws.Range("B318").Value = 63565
$ws.Range("C318").Value = "HUL-Bru Inst Poly 50g"
$ws.Range("D318").Value = 102.71
$ws.Range("E318").Value = 109.19
$ws.Range("F318").Value = 60
$ws.Range("G318").Value = 6162.6

# Row 346
$ws.Range("B346").Value = 55373
$ws.Range("C346").Value = "HUL-Kissan nango jam 490g"
$ws.Range("D346").Value = 144.28
$ws.Range("E346").Value = 163.62
$ws.Range("F346").Value = -94
$ws.Range("G346").Value = -13562.32

# Row 347
$ws.Range("B347").Value = 63520
$ws.Range("C347").Value = "HUL-Kissan nango jam 490g"
$ws.Range("D347").Value = 144.28
$ws.Range("E347").Value = 153.4
$ws.Range("F347").Value = 97
$ws.Range("G347").Value = 13995.16

# Row 372
$ws.Range("B372").Value = 57885
$ws.Range("C372").Value = "HUL-Liril Soap 125 G"
$ws.Range("D372").Value = 52.13
$ws.Range("E372").Value = 62.28
$ws.Range("F372").Value = 4
$ws.Range("G372").Value = 208.52

# Row 373
$ws.Range("B373").Value = 63652
$ws.Range("C373").Value = "HUL-Liril Soap 125 G"
$ws.Range("D373").Value = 52.13
$ws.Range("E373").Value = 55.42
$ws.Range("F373").Value = 250
$ws.Range("G373").Value = 13032.5

# Row 375
$ws.Range("B375").Value = 61605
$ws.Range("C375").Value = "HUL-lux advanced eventoned glow 4x100"
$ws.Range("D375").Value = 111.96
$ws.Range("E375").Value = 133.78
$ws.Range("F375").Value = -13
$ws.Range("G375").Value = -1455.48

# Row 376
$ws.Range("B376").Value = 63563
$ws.Range("C376").Value = "HUL-lux advanced eventoned glow 4x100"
$ws.Range("D376").Value = 111.96
$ws.Range("E376").Value = 119.04
$ws.Range("F376").Value = 15
$ws.Range("G376").Value = 1679.4

# Row 382
$ws.Range("B382").Value = 63560
$ws.Range("C382").Value = "Hul-pears pure and gentle 3x125 gm"
$ws.Range("D382").Value = 126.86
$ws.Range("E382").Value = 134.87
$ws.Range("F382").Value = 104
$ws.Range("G382").Value = 13193.44

# Row 383
$ws.Range("B383").Value = 60325
$ws.Range("C383").Value = "Hul-pears pure and gentle 3x125 gm"
$ws.Range("D383").Value = 126.86
$ws.Range("E383").Value = 151.57
$ws.Range("F383").Value = -102
$ws.Range("G383").Value = -12939.72

# Row 400
$ws.Range("B400").Value = 62933
$ws.Range("C400").Value = "HUL-Sfxl Ew Bale 500G"
$ws.Range("D400").Value = 59.13
$ws.Range("E400").Value = 70.65000000000001
$ws.Range("F400").Value = 146
$ws.Range("G400").Value = 8632.98

# Row 401
$ws.Range("B401").Value = 57835
$ws.Range("C401").Value = "HUL-Sfxl Ew Bale 500G"
$ws.Range("D401").Value = 59.13
$ws.Range("E401").Value = 70.65000000000001
$ws.Range("F401").Value = 1
$ws.Range("G401").Value = 59.13

# Row 419
$ws.Range("B419").Value = 57856
$ws.Range("C419").Value = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Range("D419").Value = 171.33
$ws.Range("E419").Value = 204.69
$ws.Range("F419").Value = 2
$ws.Range("G419").Value = 342.66

# Row 420
$ws.Range("B420").Value = 63007
$ws.Range("C420").Value = "HUL-Surf Exl Mtc Liq Fl 1 Ltr Cp"
$ws.Range("D420").Value = 171.33
$ws.Range("E420").Value = 204.69
$ws.Range("F420").Value = 984
$ws.Range("G420").Value = 168588.72

# Row 431
$ws.Range("B431").Value = 53082
$ws.Range("C431").Value = "HUL-VIM BAR MULTIPACK FW 4X200G"
$ws.Range("D431").Value = 59.47
$ws.Range("E431").Value = 71.05
$ws.Range("F431").Value = 1
$ws.Range("G431").Value = 59.47

# Row 432
$ws.Range("B432").Value = 63102
$ws.Range("C432").Value = "HUL-Vim Bar Multipack Fw 4X200G"
$ws.Range("D432").Value = 59.47
$ws.Range("E432").Value = 71.05
$ws.Range("F432").Value = 36
$ws.Range("G432").Value = 2140.92

# Row 536
$ws.Range("B536").Value = 47097
$ws.Range("C536").Value = "KUS-Floor Wiper"
$ws.Range("D536").Value = 112.28
$ws.Range("E536").Value = 134.16
$ws.Range("F536").Value = 15
$ws.Range("G536").Value = 1684.2

# Row 537
$ws.Range("B537").Value = 58047
$ws.Range("C537").Value = "KUS-Floor Wiper"
$ws.Range("D537").Value = 105.54
$ws.Range("E537").Value = 126.1
$ws.Range("F537").Value = 54
$ws.Range("G537").Value = 5699.16

# Row 579
$ws.Range("B579").Value = 65069
$ws.Range("C579").Value = "CRE-Bourbon 100gm"
$ws.Range("D579").Value = 13.45
$ws.Range("E579").Value = 14.3
$ws.Range("F579").Value = 172
$ws.Range("G579").Value = 2313.4

# Row 580
$ws.Range("B580").Value = 53757
$ws.Range("C580").Value = "CRE-Bourbon 100gm"
$ws.Range("D580").Value = 13.45
$ws.Range("E580").Value = 16.08
$ws.Range("F580").Value = -159
$ws.Range("G580").Value = -2138.55

# Row 590
$ws.Range("B590").Value = 45706
$ws.Range("C590").Value = "CRE-Cremica Golden Bytes Rich Butter 200Gm"
$ws.Range("D590").Value = 19.73
$ws.Range("E590").Value = 23.58
$ws.Range("F590").Value = -202
$ws.Range("G590").Value = -3985.46

# Row 591
$ws.Range("B591").Value = 64922
$ws.Range("C591").Value = "CRE-Cremica Golden Bytes Rich Butter 200Gm"
$ws.Range("D591").Value = 19.73
$ws.Range("E591").Value = 20.98
$ws.Range("F591").Value = 207
$ws.Range("G591").Value = 4084.11

# Row 599
$ws.Range("B599").Value = 64925
$ws.Range("C599").Value = "CRE-Cremica Oatmeal Digestive 112.5 Gm"
$ws.Range("D599").Value = 13.15
$ws.Range("E599").Value = 13.97
$ws.Range("F599").Value = 302
$ws.Range("G599").Value = 3971.3

# Row 600
$ws.Range("B600").Value = 45709
$ws.Range("C600").Value = "CRE-Cremica Oatmeal Digestive 112.5 Gm"
$ws.Range("D600").Value = 13.15
$ws.Range("E600").Value = 15.69
$ws.Range("F600").Value = -300
$ws.Range("G600").Value = -3945

# Row 601
$ws.Range("B601").Value = 64919
$ws.Range("C601").Value = "CRE-Cremica Pista Almond Cookies (75 +25Gm)"
$ws.Range("D601").Value = 26.3
$ws.Range("E601").Value = 27.97
$ws.Range("F601").Value = 224
$ws.Range("G601").Value = 5891.2

# Row 602
$ws.Range("B602").Value = 45702
$ws.Range("C602").Value = "CRE-Cremica Pista Almond Cookies (75 +25Gm)"
$ws.Range("D602").Value = 26.3
$ws.Range("E602").Value = 31.43
$ws.Range("F602").Value = -215
$ws.Range("G602").Value = -5654.5

# Row 604
$ws.Range("B604").Value = 65067
$ws.Range("C604").Value = "CRE-Kaju khz cookies 100 gm"
$ws.Range("D604").Value = 14.73
$ws.Range("E604").Value = 15.65
$ws.Range("F604").Value = 338
$ws.Range("G604").Value = 4978.74

# Row 605
$ws.Range("B605").Value = 53595
$ws.Range("C605").Value = "CRE-Kaju khz cookies 100 gm"
$ws.Range("D605").Value = 14.73
$ws.Range("E605").Value = 17.61
$ws.Range("F605").Value = -335
$ws.Range("G605").Value = -4934.55
